$wb = $excel.ActiveWorkbook

# TEST_CASES is the first sheet in the workbook.
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Data change -----------------------------------------------------
# Row 2 held a scripted test case (TC_KIND = SCRIPTED). It is being
# changed to a "TRAP !!!" kind so the TC_SCRIPTING_LANGUAGE /
# TC_SCRIPT columns now carry language info for the research screens.
# TC_KIND (X2) switches from "SCRIPTED" to "TRAP !!!"; the previously
# unused "SCRIPTED" shared string is dropped automatically once nothing
# references it any more. TC_SCRIPTING_LANGUAGE (Y2) and TC_SCRIPT (Z2)
# keep their existing text values.
$ws.Range("X2").Value = "TRAP !!!"
$ws.Range("Y2").Value = "TRAP !!!"
$ws.Range("Z2").Value = "Feature: Make something"

# --- View state --------------------------------------------------------
# The sheet's scroll position / selection moved from K1 / W10 to N1 / X5.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 14
$ws.Range("X5").Select()
